$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 400
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -87
$ws.Range("H11").Value = 273.6
$ws.Range("I11").Value = 273.6
$ws.Range("K11").Value = 273.6
$ws.Range("M11").Value = -133.6
$ws.Range("H28").Value = 655.3461
$ws.Range("I28").Value = 366.38095
$ws.Range("K28").Value = 366.38095
$ws.Range("M28").Value = 118.61905
$ws.Range("H92").Value = 3115.0527
$ws.Range("I92").Value = 6399.8335
$ws.Range("J92").Value = 1599
$ws.Range("K92").Value = 6399.8335
$ws.Range("L92").Value = 1599
$ws.Range("M92").Value = -5151.8335
$ws.Range("N92").Value = -4095
$ws.Range("H125").Value = 4136842.8
$ws.Range("I125").Value = 6499225
$ws.Range("J125").Value = 2673.5
$ws.Range("K125").Value = 58493025
$ws.Range("L125").Value = 24061.5
$ws.Range("M125").Value = -58490565
$ws.Range("N125").Value = -28981.5
$ws.Range("H129").Value = 1032
$ws.Range("I129").Value = 584.7368
$ws.Range("J129").Value = 3156.5
$ws.Range("K129").Value = 1754.2104
$ws.Range("L129").Value = 9469.5
$ws.Range("M129").Value = 3245.7896
$ws.Range("N129").Value = -19469.5
$ws.Range("H137").Value = 15496.75
$ws.Range("I137").Value = 1631.091
$ws.Range("J137").Value = 46001.2
$ws.Range("K137").Value = 4893.272999999999
$ws.Range("L137").Value = 138003.6
$ws.Range("M137").Value = -2343.272999999999
$ws.Range("N137").Value = -143103.6
$ws.Range("H138").Value = 2473.0847
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2473.0847
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7419.2541
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -17699.2541

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1255669.8
$ws.Range("I5").Value = 3390107
$ws.Range("K5").Value = 3390107
$ws.Range("M5").Value = -3389995
$ws.Range("H88").Value = 1852.4546
$ws.Range("J88").Value = 2109.2856
$ws.Range("L88").Value = 2109.2856
$ws.Range("N88").Value = -2921.2856
$ws.Range("H91").Value = 1852.4546
$ws.Range("J91").Value = 2109.2856
$ws.Range("L91").Value = 2109.2856
$ws.Range("N91").Value = -4917.2856
$ws.Range("H122").Value = 1153921.2
$ws.Range("I122").Value = 1803678.1
$ws.Range("J122").Value = 4351.3076
$ws.Range("K122").Value = 5411034.300000001
$ws.Range("L122").Value = 13053.9228
$ws.Range("M122").Value = -5408584.300000001
$ws.Range("N122").Value = -17953.9228
$ws.Range("H132").Value = 4563704.5
$ws.Range("J132").Value = 10033384
$ws.Range("L132").Value = 30100152
$ws.Range("N132").Value = -30105212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1255669.8
$ws.Range("I4").Value = 3390107
$ws.Range("K4").Value = 3390107
$ws.Range("M4").Value = -3389992
$ws.Range("H134").Value = 28012.725
$ws.Range("I134").Value = 29606.389
$ws.Range("K134").Value = 88819.167
$ws.Range("M134").Value = -86284.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27215
$ws.Range("J31").Value = 45318.855
$ws.Range("L31").Value = 45318.855
$ws.Range("N31").Value = -45908.855
$ws.Range("H34").Value = 27215
$ws.Range("J34").Value = 45318.855
$ws.Range("L34").Value = 45318.855
$ws.Range("N34").Value = -45722.855
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3774
$ws.Range("I3").Value = 3774
$ws.Range("K3").Value = 11322
$ws.Range("M3").Value = -11210
$ws.Range("H133").Value = 6112.5713
$ws.Range("I133").Value = 3997
$ws.Range("J133").Value = 8933.333000000001
$ws.Range("K133").Value = 11991
$ws.Range("L133").Value = 26799.999
$ws.Range("M133").Value = -6931
$ws.Range("N133").Value = -36919.999
$ws.Range("H134").Value = 5780.1353
$ws.Range("I134").Value = 1441.6666
$ws.Range("J134").Value = 6619.839
$ws.Range("K134").Value = 4324.9998
$ws.Range("L134").Value = 19859.517
$ws.Range("M134").Value = 745.0002000000004
$ws.Range("N134").Value = -29999.517
$ws.Range("H137").Value = 3560.9333
$ws.Range("I137").Value = 3497.4
$ws.Range("J137").Value = 3592.7
$ws.Range("K137").Value = 10492.2
$ws.Range("L137").Value = 10778.1
$ws.Range("M137").Value = -5392.200000000001
$ws.Range("N137").Value = -20978.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14495.723
$ws.Range("J80").Value = 21999.4
$ws.Range("L80").Value = 21999.4
$ws.Range("N80").Value = -23995.4
$ws.Range("H83").Value = 14495.723
$ws.Range("J83").Value = 21999.4
$ws.Range("L83").Value = 109997
$ws.Range("N83").Value = -119981
$ws.Range("H102").Value = 5877596.5
$ws.Range("I102").Value = 13515932
$ws.Range("K102").Value = 13515932
$ws.Range("M102").Value = -13514310
$ws.Range("H123").Value = 54394.5
$ws.Range("J123").Value = 54394.5
$ws.Range("L123").Value = 54394.5
$ws.Range("N123").Value = -59294.5
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1708202.9
$ws.Range("I7").Value = 2983294.2
$ws.Range("J7").Value = 8080.9165
$ws.Range("K7").Value = 2983294.2
$ws.Range("L7").Value = 8080.9165
$ws.Range("M7").Value = -2983182.2
$ws.Range("N7").Value = -8304.916499999999
$ws.Range("H22").Value = 47621104
$ws.Range("I22").Value = 1770
$ws.Range("J22").Value = 90911410
$ws.Range("K22").Value = 1770
$ws.Range("L22").Value = 90911410
$ws.Range("M22").Value = -1475
$ws.Range("N22").Value = -90912000
$ws.Range("H27").Value = 47621104
$ws.Range("I27").Value = 1770
$ws.Range("J27").Value = 90911410
$ws.Range("K27").Value = 1770
$ws.Range("L27").Value = 90911410
$ws.Range("M27").Value = -1663
$ws.Range("N27").Value = -90911624
$ws.Range("H82").Value = 1852.4445
$ws.Range("I82").Value = 2287.15
$ws.Range("J82").Value = 1309.0625
$ws.Range("K82").Value = 2287.15
$ws.Range("L82").Value = 1309.0625
$ws.Range("M82").Value = -1926.15
$ws.Range("N82").Value = -2031.0625
$ws.Range("H85").Value = 1852.4445
$ws.Range("I85").Value = 2287.15
$ws.Range("J85").Value = 1309.0625
$ws.Range("K85").Value = 2287.15
$ws.Range("L85").Value = 1309.0625
$ws.Range("M85").Value = -1039.15
$ws.Range("N85").Value = -3805.0625
$ws.Range("H126").Value = 1708202.9
$ws.Range("I126").Value = 2983294.2
$ws.Range("J126").Value = 8080.9165
$ws.Range("K126").Value = 8949882.600000001
$ws.Range("L126").Value = 24242.7495
$ws.Range("M126").Value = -8947412.600000001
$ws.Range("N126").Value = -29182.7495
$ws.Range("H136").Value = 10110.143
$ws.Range("I136").Value = 7491.6284
$ws.Range("J136").Value = 16656.428
$ws.Range("K136").Value = 22474.8852
$ws.Range("L136").Value = 49969.284
$ws.Range("M136").Value = -19924.8852
$ws.Range("N136").Value = -55069.284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 639271.6
$ws.Range("I122").Value = 1012836.25
$ws.Range("J122").Value = 7085.385
$ws.Range("K122").Value = 3038508.75
$ws.Range("L122").Value = 21256.155
$ws.Range("M122").Value = -3036058.75
$ws.Range("N122").Value = -26156.155
